# Add the newly-researched "十六国" (Sixteen Kingdoms) + Three-Kingdoms era
# rows to the "年号" (era-name) summary sheet, matching the source data
# already present on the other sheets of this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "年号" sheet

$data = @(
  ,@("汉赵", "元熙", "刘渊", 304, 308)
  ,@("汉赵", "永凤", "刘渊", 308, 309)
  ,@("汉赵", "河瑞", "刘渊", 309, 310)
  ,@("汉赵", "光兴", "刘聪", 310, 311)
  ,@("汉赵", "嘉平", "刘聪", 311, 315)
  ,@("汉赵", "建元", "刘聪", 315, 316)
  ,@("汉赵", "鹏嘉", "刘聪", 316, 318)
  ,@("汉赵", "汉昌", "刘粲", 318, 318)
  ,@("汉赵", "光初", "刘曜", 318, 329)
  ,@("成汉", "建初", "李特", 303, 304)
  ,@("成汉", "建兴", "李雄", 304, 306)
  ,@("成汉", "昭平", "李雄", 306, 310)
  ,@("成汉", "玉衡", "李雄", 311, 334)
  ,@("成汉", "玉恒", "李期", 335, 338)
  ,@("成汉", "汉兴", "李寿", 338, 343)
  ,@("成汉", "大和", "李寿", 344, 346)
  ,@("成汉", "嘉宁", "李寿", 346, 347)
  ,@("后赵", "太和", "石勒", 328, 330)
  ,@("后赵", "建平", "石勒", 330, 333)
  ,@("后赵", "延熙", "石弘", 334, 334)
  ,@("后赵", "建武", "石虎", 334, 349)
  ,@("后赵", "太宁", "石遵", 349, 350)
  ,@("后赵", "青龙", "石祗", 350, 350)
  ,@("后赵", "永兴", "石祗", 350, 351)
  ,@("前燕", "燕元", "慕容皝", 345, 348)
  ,@("前燕", "元玺", "慕容俊", 352, 357)
  ,@("前燕", "光寿", "慕容俊", 357, 359)
  ,@("前燕", "建熙", "慕容暐", 360, 370)
  ,@("前秦", "皇始", "苻健", 351, 355)
  ,@("前秦", "寿光", "苻生", 355, 357)
  ,@("前秦", "甘露", "苻坚", 359, 364)
  ,@("前秦", "建元", "苻坚", 365, 385)
  ,@("前秦", "太安", "苻丕", 385, 386)
  ,@("前秦", "弘始", "苻登", 386, 394)
  ,@("前秦", "永安", "苻崇", 394, 394)
  ,@("曹魏", "黄初", "曹丕", 220, 226)
  ,@("曹魏", "太和", "曹叡", 227, 233)
  ,@("曹魏", "青龙", "曹叡", 233, 237)
  ,@("曹魏", "景初", "曹叡", 237, 239)
  ,@("曹魏", "正始", "曹芳", 240, 249)
  ,@("曹魏", "嘉平", "曹髦", 249, 254)
  ,@("曹魏", "正元", "曹髦", 254, 256)
  ,@("曹魏", "甘露", "曹髦", 256, 260)
  ,@("曹魏", "景元", "曹奂", 260, 264)
  ,@("曹魏", "咸熙", "曹奂", 264, 265)
  ,@("蜀汉", "章武", "刘备", 221, 223)
  ,@("蜀汉", "建兴", "刘禅", 223, 237)
  ,@("蜀汉", "延熙", "刘禅", 238, 257)
  ,@("蜀汉", "景耀", "刘禅", 258, 263)
  ,@("蜀汉", "炎兴", "刘禅", 263, 263)
  ,@("孙吴", "黄武", "孙权", 222, 229)
  ,@("孙吴", "黄龙", "孙权", 229, 231)
  ,@("孙吴", "嘉禾", "孙权", 232, 238)
  ,@("孙吴", "赤乌", "孙权", 238, 251)
  ,@("孙吴", "太元", "孙权", 251, 252)
  ,@("孙吴", "神凤", "孙亮", 252, 252)
  ,@("孙吴", "建兴", "孙亮", 252, 253)
  ,@("孙吴", "五凤", "孙亮", 254, 256)
  ,@("孙吴", "太平", "孙休", 256, 258)
  ,@("孙吴", "永安", "孙休", 258, 264)
  ,@("孙吴", "元兴", "孙皓", 264, 265)
  ,@("孙吴", "甘露", "孙皓", 265, 266)
  ,@("孙吴", "宝鼎", "孙皓", 266, 269)
  ,@("孙吴", "建衡", "孙皓", 269, 271)
  ,@("孙吴", "凤皇", "孙皓", 272, 274)
  ,@("孙吴", "天册", "孙皓", 275, 276)
  ,@("孙吴", "天玺", "孙皓", 276, 276)
  ,@("孙吴", "天纪", "孙皓", 277, 280)

)

$startRow = 37
$r = $startRow
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $r = $r + 1
}
$endRow = $r - 1

# Match the formatting already used elsewhere on this sheet (vertical-center +
# wrap-text) by copying the format of an existing cell that carries it, then
# pasting formats only across the whole newly-added block.
$ws.Cells.Item(16, 2).Copy()
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 5)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Bring the "年号" sheet to the front (it becomes the active tab) and leave the
# selection on the newly added 孙吴 (Eastern Wu) block, as in the saved file.
$ws.Activate()
$ws.Range("A87:A104").Select()
